# June 7th, second pass: add a new line "改的第一遍" right before the existing
# paragraph that already reads "改的第一遍", and turn that original paragraph's
# wording into "改的第二遍" (as if the author typed a new first line, pressed
# Enter, and then revised the old line).

$d = $word.ActiveDocument

# Locate the (currently unique) paragraph that holds "改的第一遍".
$targetIdx = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "改的第一遍*") {
        $targetIdx = $i
        break
    }
}

if ($targetIdx -gt 1) {
    # Put the cursor at the end of the paragraph right above the target and
    # press "Enter" there, mirroring how Word would naturally insert a new
    # paragraph (and inherit the preceding run's east-Asia formatting hint).
    $prevPara = $d.Paragraphs.Item($targetIdx - 1)
    $prevRange = $prevPara.Range
    $prevRange.Collapse(0)
    $prevRange.InsertParagraphAfter()

    # Type the first-pass text into the freshly created (now empty) paragraph.
    $newPara = $d.Paragraphs.Item($targetIdx)
    $newPara.Range.Text = "改的第一遍"

    # The original paragraph has shifted down by one; update its wording to
    # reflect the second revision pass.
    $origPara = $d.Paragraphs.Item($targetIdx + 1)
    $origRange = $origPara.Range
    $origRange.Find.Execute("改的第一遍", $false, $false, $false, $false, $false, $true, 1, $false, "改的第二遍", 2)
}
